# Update score execution values
# Sheet "category_score" (the active sheet), row 12 = "SCORE EXECUTION".
# Only the execution scores for a subset of countries actually changed;
# Puerto Rico (G12) and Guatemala (H12) keep their original 7.5 value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("category_score")

$ws.Range("B12").Value = 8      # Costa Rica
$ws.Range("C12").Value = 5.1    # Chile
$ws.Range("D12").Value = 5      # Argentina
$ws.Range("E12").Value = 7      # Trinidad & Tobago
$ws.Range("F12").Value = 6.2    # Jamaica
$ws.Range("I12").Value = 6      # Columbia

# Match the row height LibreOffice recomputed for the re-wrapped text in
# row 12 after the value edits above.
$ws.Rows.Item(12).RowHeight = 23.85

# Leave the cursor where the author's last selection landed.
$ws.Range("E13").Select()
